$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'245.89"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Value = "'23.68"
$ws.Range("D3").Style = "Normal"
$ws.Range("D4").Value = "'5.363"
$ws.Range("D4").Style = "Normal"
$ws.Range("D6").Value = "'6.480"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Value = "'3.356"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Value = "'0.8113"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Value = "'0.9215"
$ws.Range("D9").Style = "Normal"
$ws.Range("D11").Value = "'0.07412"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Value = "'0.03097"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Value = "'0.03052"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Value = "'0.09362"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Value = "'3.849"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Value = "'0.001576"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Value = "'0.04706"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Value = "'0.0006042"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Value = "'0.005930"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Value = "'0.001244"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Value = "'0.004713"
$ws.Range("D21").Style = "Normal"
$ws.Range("D23").Value = "'3.597"
$ws.Range("D23").Style = "Normal"
$ws.Range("D25").Value = "'0.3230"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Value = "'0.1329"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Value = "'0.0002654"
$ws.Range("D27").Style = "Normal"
$ws.Range("D41").Value = "'0.006343"
$ws.Range("D41").Style = "Normal"
$ws.Range("D43").Value = "'0.002711"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Value = "'0.008087"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Value = "'0.00005261"
$ws.Range("D45").Style = "Normal"
$ws.Range("D47").Value = "'0.6602"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Value = "'0.001689"
$ws.Range("D48").Style = "Normal"
